$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: "SK Kemenhumkam" -> "SK Kemenhumham"
$ws.Range("E1").Value = "SK Kemenhumham"

# Update the "Nomor" value in row 2
$ws.Range("A2").Value = 12323

# Widen column E (by 6 characters, matching the author's resize)
$ws.Columns.Item(5).ColumnWidth = 18.67

# Move the active selection from F2 to A2
$ws.Range("A2").Select()
